$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet tracks GDP figures per quarterly release ("Serie"); a new
# release column for Agosto.2021 is added after the last column (BG,
# "Mayo.2021"). The new column repeats the same figures as the prior
# release (BG) until updated values are published, matching the
# "Actualización desde MV -datos-" refresh pattern seen in the rest of
# the sheet.
$ws.Range("BG1:BG19").Copy($ws.Range("BH1:BH19"))
$ws.Range("BH1").Value = "Agosto.2021"
